$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-07 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-08 Saturday", 2) | Out-Null
$d.Content.Find.Execute("99-44=", $true, $false, $false, $false, $false, $true, 1, $false, "11+17=", 2) | Out-Null
$d.Content.Find.Execute("15+8=", $true, $false, $false, $false, $false, $true, 1, $false, "33+25=", 2) | Out-Null
$d.Content.Find.Execute("56-19=", $true, $false, $false, $false, $false, $true, 1, $false, "31-1=", 2) | Out-Null
$d.Content.Find.Execute("42+7=", $true, $false, $false, $false, $false, $true, 1, $false, "67-24=", 2) | Out-Null
$d.Content.Find.Execute("39-11=", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=", 2) | Out-Null
$d.Content.Find.Execute("15+30=", $true, $false, $false, $false, $false, $true, 1, $false, "40+13=", 2) | Out-Null
$d.Content.Find.Execute("70+8=", $true, $false, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("63-18=", $true, $false, $false, $false, $false, $true, 1, $false, "23+15=", 2) | Out-Null
$d.Content.Find.Execute("63-25=", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=", 2) | Out-Null
$d.Content.Find.Execute("2+89=", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("61-22=", $true, $false, $false, $false, $false, $true, 1, $false, "42-22=", 2) | Out-Null
$d.Content.Find.Execute("23+14=", $true, $false, $false, $false, $false, $true, 1, $false, "14+37=", 2) | Out-Null
$d.Content.Find.Execute("88-25=", $true, $false, $false, $false, $false, $true, 1, $false, "13+76=", 2) | Out-Null
$d.Content.Find.Execute("39+1=", $true, $false, $false, $false, $false, $true, 1, $false, "15+73=", 2) | Out-Null
$d.Content.Find.Execute("43+25=", $true, $false, $false, $false, $false, $true, 1, $false, "70-44=", 2) | Out-Null
$d.Content.Find.Execute("46-11=", $true, $false, $false, $false, $false, $true, 1, $false, "20+70=", 2) | Out-Null
$d.Content.Find.Execute("20+0=", $true, $false, $false, $false, $false, $true, 1, $false, "67+12=", 2) | Out-Null
$d.Content.Find.Execute("22+34=", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=", 2) | Out-Null
$d.Content.Find.Execute("43-4=", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=", 2) | Out-Null
$d.Content.Find.Execute("70-50=", $true, $false, $false, $false, $false, $true, 1, $false, "69-7=", 2) | Out-Null
$d.Content.Find.Execute("69-24=", $true, $false, $false, $false, $false, $true, 1, $false, "24+66=", 2) | Out-Null
$d.Content.Find.Execute("82-50=", $true, $false, $false, $false, $false, $true, 1, $false, "55-45=", 2) | Out-Null
$d.Content.Find.Execute("37+55=", $true, $false, $false, $false, $false, $true, 1, $false, "66+4=", 2) | Out-Null
$d.Content.Find.Execute("34+3=", $true, $false, $false, $false, $false, $true, 1, $false, "49+31=", 2) | Out-Null
$d.Content.Find.Execute("24-10=", $true, $false, $false, $false, $false, $true, 1, $false, "16+24=", 2) | Out-Null
$d.Content.Find.Execute("20+4=", $true, $false, $false, $false, $false, $true, 1, $false, "54+18=", 2) | Out-Null
$d.Content.Find.Execute("63-35=", $true, $false, $false, $false, $false, $true, 1, $false, "93-0=", 2) | Out-Null
$d.Content.Find.Execute("87-70=", $true, $false, $false, $false, $false, $true, 1, $false, "75+18=", 2) | Out-Null
$d.Content.Find.Execute("30-5=", $true, $false, $false, $false, $false, $true, 1, $false, "68-2=", 2) | Out-Null
$d.Content.Find.Execute("67-13=", $true, $false, $false, $false, $false, $true, 1, $false, "27+38=", 2) | Out-Null
$d.Content.Find.Execute("83+0=", $true, $false, $false, $false, $false, $true, 1, $false, "85-63=", 2) | Out-Null
$d.Content.Find.Execute("70-48=", $true, $false, $false, $false, $false, $true, 1, $false, "73+24=", 2) | Out-Null
$d.Content.Find.Execute("62+31=", $true, $false, $false, $false, $false, $true, 1, $false, "56+19=", 2) | Out-Null
$d.Content.Find.Execute("10+52=", $true, $false, $false, $false, $false, $true, 1, $false, "58-42=", 2) | Out-Null
$d.Content.Find.Execute("48+12=", $true, $false, $false, $false, $false, $true, 1, $false, "77-74=", 2) | Out-Null
$d.Content.Find.Execute("53-2=", $true, $false, $false, $false, $false, $true, 1, $false, "42+13=", 2) | Out-Null
$d.Content.Find.Execute("14-9=", $true, $false, $false, $false, $false, $true, 1, $false, "0+89=", 2) | Out-Null
$d.Content.Find.Execute("37+57=", $true, $false, $false, $false, $false, $true, 1, $false, "81-76=", 2) | Out-Null
$d.Content.Find.Execute("96-50=", $true, $false, $false, $false, $false, $true, 1, $false, "59-43=", 2) | Out-Null
$d.Content.Find.Execute("32+26=", $true, $false, $false, $false, $false, $true, 1, $false, "92-4=", 2) | Out-Null
$d.Content.Find.Execute("34-30=", $true, $false, $false, $false, $false, $true, 1, $false, "67-27=", 2) | Out-Null
$d.Content.Find.Execute("98-3=", $true, $false, $false, $false, $false, $true, 1, $false, "43+12=", 2) | Out-Null
$d.Content.Find.Execute("94-89=", $true, $false, $false, $false, $false, $true, 1, $false, "40+22=", 2) | Out-Null
$d.Content.Find.Execute("86+6=", $true, $false, $false, $false, $false, $true, 1, $false, "44+29=", 2) | Out-Null
$d.Content.Find.Execute("55+43=", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=", 2) | Out-Null
$d.Content.Find.Execute("39+43=", $true, $false, $false, $false, $false, $true, 1, $false, "74+21=", 2) | Out-Null
$d.Content.Find.Execute("78-17=", $true, $false, $false, $false, $false, $true, 1, $false, "19+62=", 2) | Out-Null
$d.Content.Find.Execute("85-51=", $true, $false, $false, $false, $false, $true, 1, $false, "69+6=", 2) | Out-Null
$d.Content.Find.Execute("75+12=", $true, $false, $false, $false, $false, $true, 1, $false, "8+29=", 2) | Out-Null
$d.Content.Find.Execute("76-68=", $true, $false, $false, $false, $false, $true, 1, $false, "31-5=", 2) | Out-Null
$d.Content.Find.Execute("14+6=", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=", 2) | Out-Null
$d.Content.Find.Execute("55+0=", $true, $false, $false, $false, $false, $true, 1, $false, "14+2=", 2) | Out-Null
$d.Content.Find.Execute("7+75=", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=", 2) | Out-Null
$d.Content.Find.Execute("94-11=", $true, $false, $false, $false, $false, $true, 1, $false, "52-41=", 2) | Out-Null
$d.Content.Find.Execute("0+15=", $true, $false, $false, $false, $false, $true, 1, $false, "77+15=", 2) | Out-Null
$d.Content.Find.Execute("34-26=", $true, $false, $false, $false, $false, $true, 1, $false, "91-29=", 2) | Out-Null
$d.Content.Find.Execute("4+64=", $true, $false, $false, $false, $false, $true, 1, $false, "80-52=", 2) | Out-Null
$d.Content.Find.Execute("8+30=", $true, $false, $false, $false, $false, $true, 1, $false, "2+74=", 2) | Out-Null
$d.Content.Find.Execute("94-86=", $true, $false, $false, $false, $false, $true, 1, $false, "48-46=", 2) | Out-Null
$d.Content.Find.Execute("0+53=", $true, $false, $false, $false, $false, $true, 1, $false, "2+4=", 2) | Out-Null
$d.Content.Find.Execute("79+3=", $true, $false, $false, $false, $false, $true, 1, $false, "29+2=", 2) | Out-Null
$d.Content.Find.Execute("30+27=", $true, $false, $false, $false, $false, $true, 1, $false, "8+39=", 2) | Out-Null
$d.Content.Find.Execute("63-42=", $true, $false, $false, $false, $false, $true, 1, $false, "81-8=", 2) | Out-Null
$d.Content.Find.Execute("96-19=", $true, $false, $false, $false, $false, $true, 1, $false, "34-12=", 2) | Out-Null
$d.Content.Find.Execute("70-24=", $true, $false, $false, $false, $false, $true, 1, $false, "80-25=", 2) | Out-Null
$d.Content.Find.Execute("78-34=", $true, $false, $false, $false, $false, $true, 1, $false, "34+54=", 2) | Out-Null
$d.Content.Find.Execute("34+50=", $true, $false, $false, $false, $false, $true, 1, $false, "41-31=", 2) | Out-Null
$d.Content.Find.Execute("65-22=", $true, $false, $false, $false, $false, $true, 1, $false, "45-18=", 2) | Out-Null
$d.Content.Find.Execute("44+46=", $true, $false, $false, $false, $false, $true, 1, $false, "65+9=", 2) | Out-Null
$d.Content.Find.Execute("94-63=", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=", 2) | Out-Null
$d.Content.Find.Execute("78-52=", $true, $false, $false, $false, $false, $true, 1, $false, "74-67=", 2) | Out-Null
$d.Content.Find.Execute("72-32=", $true, $false, $false, $false, $false, $true, 1, $false, "65-5=", 2) | Out-Null
$d.Content.Find.Execute("3-1=", $true, $false, $false, $false, $false, $true, 1, $false, "58-51=", 2) | Out-Null
$d.Content.Find.Execute("17+71=", $true, $false, $false, $false, $false, $true, 1, $false, "65-6=", 2) | Out-Null
$d.Content.Find.Execute("7+41=", $true, $false, $false, $false, $false, $true, 1, $false, "50-35=", 2) | Out-Null
$d.Content.Find.Execute("63-54=", $true, $false, $false, $false, $false, $true, 1, $false, "16+8=", 2) | Out-Null
$d.Content.Find.Execute("52+36=", $true, $false, $false, $false, $false, $true, 1, $false, "3+47=", 2) | Out-Null
$d.Content.Find.Execute("57-51=", $true, $false, $false, $false, $false, $true, 1, $false, "14+24=", 2) | Out-Null
$d.Content.Find.Execute("81+12=", $true, $false, $false, $false, $false, $true, 1, $false, "43+26=", 2) | Out-Null
$d.Content.Find.Execute("6-1=", $true, $false, $false, $false, $false, $true, 1, $false, "53+13=", 2) | Out-Null
$d.Content.Find.Execute("98-82=", $true, $false, $false, $false, $false, $true, 1, $false, "74-27=", 2) | Out-Null
$d.Content.Find.Execute("90-50=", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=", 2) | Out-Null
$d.Content.Find.Execute("29-10=", $true, $false, $false, $false, $false, $true, 1, $false, "28+55=", 2) | Out-Null
$d.Content.Find.Execute("96-51=", $true, $false, $false, $false, $false, $true, 1, $false, "5+37=", 2) | Out-Null
$d.Content.Find.Execute("91-6=", $true, $false, $false, $false, $false, $true, 1, $false, "93-52=", 2) | Out-Null
$d.Content.Find.Execute("69+0=", $true, $false, $false, $false, $false, $true, 1, $false, "19+30=", 2) | Out-Null
$d.Content.Find.Execute("54-26=", $true, $false, $false, $false, $false, $true, 1, $false, "22+30=", 2) | Out-Null
$d.Content.Find.Execute("37-4=", $true, $false, $false, $false, $false, $true, 1, $false, "61-53=", 2) | Out-Null
$d.Content.Find.Execute("80-74=", $true, $false, $false, $false, $false, $true, 1, $false, "74-62=", 2) | Out-Null
$d.Content.Find.Execute("60-21=", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=", 2) | Out-Null
$d.Content.Find.Execute("31+51=", $true, $false, $false, $false, $false, $true, 1, $false, "78+18=", 2) | Out-Null
$d.Content.Find.Execute("27-11=", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=", 2) | Out-Null
$d.Content.Find.Execute("59-49=", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=", 2) | Out-Null
$d.Content.Find.Execute("60-1=", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=", 2) | Out-Null
$d.Content.Find.Execute("81-79=", $true, $false, $false, $false, $false, $true, 1, $false, "61+6=", 2) | Out-Null
$d.Content.Find.Execute("93-3=", $true, $false, $false, $false, $false, $true, 1, $false, "19-18=", 2) | Out-Null
$d.Content.Find.Execute("81-14=", $true, $false, $false, $false, $false, $true, 1, $false, "63+28=", 2) | Out-Null
$d.Content.Find.Execute("18+57=", $true, $false, $false, $false, $false, $true, 1, $false, "49-47=", 2) | Out-Null
$d.Content.Find.Execute("1+68=", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=", 2) | Out-Null
